# Apply updated schedule-scrape data (refresh run at 11:13:01) to all three sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912 ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 11:13:01"
$ws1.Cells.Item(3,1).Value = "Total filas: 145"
$ws1.Cells.Item(28,1).Value = "07:15:48"
$ws1.Cells.Item(28,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28,4).Value = 6
$ws1.Cells.Item(29,1).Value = "06:56:24"
$ws1.Cells.Item(29,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(29,4).Value = 25
$ws1.Cells.Item(41,1).Value = "06:38:54"
$ws1.Cells.Item(41,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(41,4).Value = 82
$ws1.Cells.Item(42,1).Value = "07:52:32"
$ws1.Cells.Item(42,3).Value = "17_ROMERO"
$ws1.Cells.Item(42,4).Value = 8
$ws1.Cells.Item(75,1).Value = "08:30:14"
$ws1.Cells.Item(75,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(75,4).Value = 47
$ws1.Cells.Item(76,1).Value = "08:52:33"
$ws1.Cells.Item(76,3).Value = "14_ABASTO"
$ws1.Cells.Item(76,4).Value = 25
$ws1.Cells.Item(77,1).Value = "08:40:59"
$ws1.Cells.Item(77,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(77,4).Value = 37
$ws1.Cells.Item(78,1).Value = "08:30:14"
$ws1.Cells.Item(78,3).Value = "14_ABASTO"
$ws1.Cells.Item(78,4).Value = 48
$ws1.Cells.Item(79,1).Value = "08:52:33"
$ws1.Cells.Item(79,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(79,4).Value = 26
$ws1.Cells.Item(85,1).Value = "08:30:14"
$ws1.Cells.Item(85,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(85,4).Value = 72
$ws1.Cells.Item(86,1).Value = "08:40:59"
$ws1.Cells.Item(86,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(86,4).Value = 62
$ws1.Cells.Item(118,1).Value = "11:13:01"
$ws1.Cells.Item(118,4).Value = 2
$ws1.Cells.Item(119,1).Value = "11:13:01"
$ws1.Cells.Item(119,4).Value = 2
$ws1.Cells.Item(122,1).Value = "11:13:01"
$ws1.Cells.Item(122,4).Value = 12
$ws1.Cells.Item(123,1).Value = "11:13:01"
$ws1.Cells.Item(123,4).Value = 16
$ws1.Cells.Item(125,1).Value = "11:13:01"
$ws1.Cells.Item(125,4).Value = 18
$ws1.Cells.Item(127,1).Value = "11:13:01"
$ws1.Cells.Item(127,4).Value = 29
$ws1.Cells.Item(128,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(129,1).Value = "11:13:01"
$ws1.Cells.Item(129,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(129,4).Value = 32
$ws1.Cells.Item(131,1).Value = "11:13:01"
$ws1.Cells.Item(131,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(131,4).Value = 39
$ws1.Cells.Item(132,1).Value = "10:07:51"
$ws1.Cells.Item(132,2).Value = "11:52"
$ws1.Cells.Item(132,4).Value = 105
$ws1.Cells.Item(133,1).Value = "11:13:01"
$ws1.Cells.Item(133,2).Value = "11:53"
$ws1.Cells.Item(133,3).Value = "225_GOMEZ"
$ws1.Cells.Item(133,4).Value = 40
$ws1.Cells.Item(134,1).Value = "11:13:01"
$ws1.Cells.Item(134,2).Value = "11:58"
$ws1.Cells.Item(134,3).Value = "17_ROMERO"
$ws1.Cells.Item(134,4).Value = 45
$ws1.Cells.Item(135,1).Value = "10:41:48"
$ws1.Cells.Item(135,2).Value = "12:05"
$ws1.Cells.Item(135,4).Value = 84
$ws1.Cells.Item(136,1).Value = "11:13:01"
$ws1.Cells.Item(136,2).Value = "12:06"
$ws1.Cells.Item(136,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(136,4).Value = 53
$ws1.Cells.Item(137,1).Value = "11:13:01"
$ws1.Cells.Item(137,4).Value = 57
$ws1.Cells.Item(138,1).Value = "11:13:01"
$ws1.Cells.Item(138,2).Value = "12:10"
$ws1.Cells.Item(138,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(138,4).Value = 57
$ws1.Cells.Item(139,1).Value = "11:13:01"
$ws1.Cells.Item(139,2).Value = "12:17"
$ws1.Cells.Item(139,3).Value = "10_OLMOS"
$ws1.Cells.Item(139,4).Value = 64
$ws1.Cells.Item(140,2).Value = "12:21"
$ws1.Cells.Item(140,3).Value = "215C_EL PATO"
$ws1.Cells.Item(140,4).Value = 100
$ws1.Cells.Item(141,1).Value = "11:13:01"
$ws1.Cells.Item(141,2).Value = "12:22"
$ws1.Cells.Item(141,3).Value = "215C_EL PATO"
$ws1.Cells.Item(141,4).Value = 69
$ws1.Cells.Item(142,1).Value = "11:13:01"
$ws1.Cells.Item(142,2).Value = "12:31"
$ws1.Cells.Item(142,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(142,4).Value = 78
$ws1.Cells.Item(143,1).Value = "10:41:48"
$ws1.Cells.Item(143,2).Value = "12:32"
$ws1.Cells.Item(143,3).Value = "14_ABASTO"
$ws1.Cells.Item(143,4).Value = 111
$ws1.Cells.Item(144,1).Value = "11:13:01"
$ws1.Cells.Item(144,2).Value = "12:33"
$ws1.Cells.Item(144,3).Value = "14_ABASTO"
$ws1.Cells.Item(144,4).Value = 80
$ws1.Cells.Item(145,1).Value = "11:13:01"
$ws1.Cells.Item(145,2).Value = "12:33"
$ws1.Cells.Item(145,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(145,4).Value = 80
$ws1.Cells.Item(146,1).Value = "11:13:01"
$ws1.Cells.Item(146,2).Value = "12:34"
$ws1.Cells.Item(146,3).Value = "15_ABASTO"
$ws1.Cells.Item(146,4).Value = 81
$ws1.Cells.Item(146,5).Value = "LP1912"
$ws1.Cells.Item(147,1).Value = "10:56:01"
$ws1.Cells.Item(147,2).Value = "12:34"
$ws1.Cells.Item(147,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(147,4).Value = 98
$ws1.Cells.Item(147,5).Value = "LP1912"
$ws1.Cells.Item(148,1).Value = "10:41:48"
$ws1.Cells.Item(148,2).Value = "12:36"
$ws1.Cells.Item(148,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(148,4).Value = 115
$ws1.Cells.Item(148,5).Value = "LP1912"
$ws1.Cells.Item(149,1).Value = "11:13:01"
$ws1.Cells.Item(149,2).Value = "12:48"
$ws1.Cells.Item(149,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(149,4).Value = 95
$ws1.Cells.Item(149,5).Value = "LP1912"
$ws1.Cells.Item(150,1).Value = "11:13:01"
$ws1.Cells.Item(150,2).Value = "13:03"
$ws1.Cells.Item(150,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(150,4).Value = 110
$ws1.Cells.Item(150,5).Value = "LP1912"

# --- LP1912-215 ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 11:13:01"
$ws2.Cells.Item(22,1).Value = "11:13:01"
$ws2.Cells.Item(22,4).Value = 18
$ws2.Cells.Item(24,1).Value = "11:13:01"
$ws2.Cells.Item(24,4).Value = 29
$ws2.Cells.Item(26,1).Value = "11:13:01"
$ws2.Cells.Item(26,4).Value = 69

# --- 6203-6173 ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 11:13:01"
$ws3.Cells.Item(23,1).Value = "11:13:01"
$ws3.Cells.Item(23,4).Value = 13
